$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in previously empty cells with new data ---

# Row 11: Apollo Intensa Emozione
$ws.Range("C11").Value = 2018
$ws.Range("D11").Value = "RWD"
$ws.Range("E11").Value = 773
$ws.Range("F11").Value = "H"
$ws.Range("G11").Value = "épique"
$ws.Range("H11").Value = 0

# Row 55
$ws.Range("D55").Value = "RWD"
$ws.Range("G55").Value = "épique"
$ws.Range("H55").Value = 0

# Row 63
$ws.Range("D63").Value = "4WD"
$ws.Range("G63").Value = "rare"
$ws.Range("H63").Value = 0

# Row 75
$ws.Range("D75").Value = "4WD"
$ws.Range("G75").Value = "rare"
$ws.Range("H75").Value = 0

# Row 80
$ws.Range("D80").Value = "RWD"
$ws.Range("G80").Value = "légendaire"
$ws.Range("H80").Value = 0

# Row 90
$ws.Range("D90").Value = "RWD"
$ws.Range("G90").Value = "rare"
$ws.Range("H90").Value = 0

# Row 101
$ws.Range("D101").Value = "RWD"
$ws.Range("G101").Value = "épique"
$ws.Range("H101").Value = 0

# --- Column G width (best-fit after data entry) ---
$ws.Columns.Item(7).ColumnWidth = 15.3

# --- Row heights recompute (wrap-text autofit after column width change) ---
$ws.Rows.Item(4).RowHeight = 13.5
$ws.Rows.Item(5).RowHeight = 13.5
$ws.Rows.Item(6).RowHeight = 18.75
$ws.Rows.Item(8).RowHeight = 12
$ws.Rows.Item(11).RowHeight = 18.75
$ws.Rows.Item(13).RowHeight = 16.5
$ws.Rows.Item(14).RowHeight = 15
$ws.Rows.Item(15).RowHeight = 18
$ws.Rows.Item(22).RowHeight = 12
$ws.Rows.Item(31).RowHeight = 14.25
$ws.Rows.Item(37).RowHeight = 16.5
$ws.Rows.Item(66).RowHeight = 18.75
$ws.Rows.Item(88).RowHeight = 15.75
$ws.Rows.Item(89).RowHeight = 19.5
$ws.Rows.Item(93).RowHeight = 18
$ws.Rows.Item(97).RowHeight = 15
$ws.Rows.Item(100).RowHeight = 14.25

# --- Selection / view state: scroll back to top, select H11 ---
$ws.Range("A1").Select() | Out-Null
$ws.Range("H11").Select() | Out-Null
